$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row (row 1) values for columns B-E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 (CON) values for columns B-E
$ws.Range("B2").Value = -8.6796474851898218
$ws.Range("C2").Value = 6.5796205393048064
$ws.Range("D2").Value = 5.1667374473021255
$ws.Range("E2").Value = 7.6813126148273794

# Update row 3 (STR) values for columns B-E
$ws.Range("B3").Value = 7.1691463110793299
$ws.Range("C3").Value = 18.788575696066857
$ws.Range("D3").Value = 30.181396220663299
$ws.Range("E3").Value = 1.1626872691771324

# Reflect the updated selection range used during editing
$ws.Range("B1:E3").Select()
